$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty cells for the 0.25 (06:00) hour row with new production case counts
$ws.Range("B4").Value = 161
$ws.Range("G4").Value = 164
$ws.Range("L4").Value = 125

# Update the active selection to match the final state (L10)
$ws.Range("L10").Select()

$wb.Save()
